# mudanca no layout e testes
# Populate the price-list sheet with the quote rows (descricao, modelo,
# fabricante, fornecedor, quantidade, valor unitario, data cotacao, total)
# and a grand-total row, alternating a light-grey / light-blue row fill on
# top of the existing blue header band.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlHAlignLeft = -4131
$xlHAlignCenter = -4108
$xlVAlignCenter = -4108

# BGR-packed ints for the COM .Color properties (0x00BBGGRR)
$cGrey     = 14540253   # DDDDDD
$cBlue     = 15853019   # dbe5f1
$cGreen    = 8388352    # 00FF7F - same green already used by the old total cell
$cItalicFg = 5263440    # 505050

$xlLineStyleNone = -4142

function Format-Cell($cell, $fillColor, $bold, $italic, $hAlign) {
    # the template's row-2 band carried a top border; strip any inherited
    # border before laying down the new (borderless) body-row style so we
    # don't drag stale borders into the new rows
    $cell.Borders.Item(8).LineStyle = $xlLineStyleNone
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.Bold = $bold
    $cell.Font.Italic = $italic
    if ($italic) {
        $cell.Font.Color = $cItalicFg
    }
    $cell.Interior.Color = $fillColor
    $cell.HorizontalAlignment = $hAlign
    $cell.VerticalAlignment = $xlVAlignCenter
}

function Format-DataRow($rowNum, $fillColor) {
    # A-D: descricao/modelo/fabricante/fornecedor -> left aligned, plain font
    Format-Cell $ws.Cells.Item($rowNum, 1) $fillColor $false $false $xlHAlignLeft
    Format-Cell $ws.Cells.Item($rowNum, 2) $fillColor $false $false $xlHAlignLeft
    Format-Cell $ws.Cells.Item($rowNum, 3) $fillColor $false $false $xlHAlignLeft
    Format-Cell $ws.Cells.Item($rowNum, 4) $fillColor $false $false $xlHAlignLeft
    # E: quantidade -> centered, plain font
    Format-Cell $ws.Cells.Item($rowNum, 5) $fillColor $false $false $xlHAlignCenter
    # F: valor unitario -> centered, italic grey font
    Format-Cell $ws.Cells.Item($rowNum, 6) $fillColor $false $true  $xlHAlignCenter
    # G: data cotacao -> centered, plain font
    Format-Cell $ws.Cells.Item($rowNum, 7) $fillColor $false $false $xlHAlignCenter
    # H: total -> centered, bold font
    Format-Cell $ws.Cells.Item($rowNum, 8) $fillColor $true  $false $xlHAlignCenter
}

$rows = @(
    @{ row=2; a="CAIXA DE PASSAGEM DE SOBREPOR 120x120x7,5CM"; b=$null;                     c="Furukawa";    d="Nucleo";    e=5;  f=10;    g="13/05/2020"; fill=$cGrey },
    @{ row=3; a="Câmera mini Bullet";                          b="DS-2CD2012-I";             c="Hikvision";   d="Hikvision"; e=12; f=329.9; g="07/05/2020"; fill=$cBlue },
    @{ row=4; a='ELETRODUTO GALVANIZADO A FOGO DE 1"';         b=$null;                     c="Ferro Norte"; d="Matec";     e=50; f=15;    g="17/07/2020"; fill=$cGrey },
    @{ row=5; a="ELETRODUTO PVC RÍGIDO DE ½” ANTICHAMA";       b=$null;                     c="Ferro Norte"; d="Matec";     e=32; f=7.5;   g="15/05/2020"; fill=$cBlue },
    @{ row=6; a="Guia de cabo";                                b=$null;                     c="furukawa";    d="Matec";     e=3;  f=25;    g="18/11/1991"; fill=$cGrey },
    @{ row=7; a="Patch Cord cat 6";                            b=$null;                     c="Furukawa";    d="Engecopi";  e=45; f=7.9;   g="01/05/2020"; fill=$cBlue },
    @{ row=8; a="Switch POE 16 Portas";                        b="DS-7716/7732NI-K4/16P";   c="Hikvision";   d="Hikvision"; e=1;  f=800;   g="15/12/2019"; fill=$cGrey }
)

foreach ($r in $rows) {
    $rowNum = $r.row
    $ws.Cells.Item($rowNum, 1).Value = $r.a
    if ($r.b) {
        $ws.Cells.Item($rowNum, 2).Value = $r.b
    }
    $ws.Cells.Item($rowNum, 3).Value = $r.c
    $ws.Cells.Item($rowNum, 4).Value = $r.d
    $ws.Cells.Item($rowNum, 5).Value = $r.e
    $ws.Cells.Item($rowNum, 6).Value = $r.f
    # leading apostrophe forces the dd/mm/yyyy text to stay a literal string
    # instead of being auto-parsed into a date serial number
    $ws.Cells.Item($rowNum, 7).Value = "'" + $r.g
    $ws.Cells.Item($rowNum, 8).Formula = "=E" + $rowNum + "*F" + $rowNum

    Format-DataRow $rowNum $r.fill
}

# ---- Row 9: grand-total row ----
for ($col = 1; $col -le 7; $col++) {
    $cell = $ws.Cells.Item(9, $col)
    $cell.Borders.Item(8).Color = 0
    $cell.Borders.Item(8).LineStyle = 1
}

$totalCell = $ws.Range("H9")
$totalCell.Formula = "=SUM(H2:H8)"
$totalCell.Font.Name = "Arial"
$totalCell.Font.Size = 10
$totalCell.Font.Bold = $true
$totalCell.Interior.Color = $cGreen
$totalCell.HorizontalAlignment = $xlHAlignCenter
$totalCell.VerticalAlignment = $xlVAlignCenter
$totalCell.Borders.Item(8).Color = 0
$totalCell.Borders.Item(8).LineStyle = 1

Write-Output "layout updated"
